# Gamma1F-HW30: two new simulation methods ("Holden" and "Rizzie Spiral")
# were inserted into the method list right after "Spiral5" (every later
# method's name/index label keeps its own row - only the C:T numeric
# columns ripple down by two rows), "Thomas Hex" was renamed to
# "Matthies Hex", and the HW30 simulation was rerun, producing new
# numbers throughout and two extra trailing rows (30 and 31) for the
# methods that no longer fit in the original 2..29 row range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# ---------------------------------------------------------------------
# 1) Ripple the C:T numeric columns down by two rows (bottom-up so a row
#    is always read before it gets overwritten). Rows 28 and 29's old
#    data ends up forming the two brand-new rows 30/31 below.
# ---------------------------------------------------------------------
for ($r = 27; $r -ge 4; $r--) {
    $dest = $r + 2
    $src = $ws.Range($cols[0] + $r + ":" + $cols[$cols.Length-1] + $r)
    $dst = $ws.Range($cols[0] + $dest + ":" + $cols[$cols.Length-1] + $dest)
    $src.Copy($dst)
}

# ---------------------------------------------------------------------
# 2) New row 30 = old row 28's data, new row 31 = old row 29's data
#    (captured into variables before step 1 overwrote rows 28/29... but
#    since step 1 only walked down to row 4->6, rows 28/29 are still the
#    original values at this point). Build the brand-new rows 30 and 31
#    (index/label columns + data) from those original values.
# ---------------------------------------------------------------------
$row30 = @(1.004366018022624, 0.9998427642715012, 1.004366018022624, 0.9982913481131679, 0.9987436622534801, 1.003042781811578, 0.9946435605007281, 0.9998427642715012, 0.9998427642715012, 0.9982913481131679, 1.001328683067896, 1.001328683067896, 1.00190004931579, 1.000833376802431, 1.000833376802431, 1.000585723669698, 1.000585723669698, 0.9998216891621797)
$row31 = @(1.10398746778623, 0.9757702402422944, 1.10398746778623, 0.9733593318942846, 0.9740622323429464, 1.056778150388713, 0.93023692249816, 0.9757702402422944, 0.9757702402422944, 0.9733593318942846, 1.038673399840257, 1.038673399840257, 1.044708316689743, 1.01770567997427, 1.01770567997427, 1.007221820041276, 1.007221820041276, 1.002365724192105)

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Michael-CCHex"
for ($i = 0; $i -lt $row30.Length; $i++) {
    $ws.Cells.Item(30, 3 + $i).Value = $row30[$i]
}

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Michael-SNHex"
for ($i = 0; $i -lt $row31.Length; $i++) {
    $ws.Cells.Item(31, 3 + $i).Value = $row31[$i]
}

# Match the header column's bordered/bold/centered label style used by
# the rest of column A.
$ws.Range("A30").Font.Bold = $true
$ws.Range("A30").HorizontalAlignment = -4108
$ws.Range("A30").VerticalAlignment = -4160
$ws.Range("A30").Borders.LineStyle = 1
$ws.Range("A31").Font.Bold = $true
$ws.Range("A31").HorizontalAlignment = -4108
$ws.Range("A31").VerticalAlignment = -4160
$ws.Range("A31").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3) Rows 4 and 5 keep their own A/B row labels (2 / Holden, 3 / Rizzie
#    Spiral - the two new methods slotted in right after "Spiral5"); set
#    their label text and freshly (re)computed C:T values.
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = "Holden"
$row4 = @(0.9283453847567166, 1.011168748892794, 0.9283453847567166, 1.021694309866923, 1.018625634874417, 0.9580576387156438, 1.060028576728711, 1.011168748892794, 1.011168748892794, 1.021694309866923, 0.9750198473118197, 0.9750198473118197, 0.969365777779761, 0.9870694811721443, 0.9870694811721442, 0.9930942981023065, 0.9930942981023065, 0.9996533823058673)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 3 + $i).Value = $row4[$i]
}

$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$row5 = @(0.9392084920655477, 1.024619394198223, 0.9392084920655477, 1.009765813781591, 1.014096298488609, 0.9705719205529828, 1.022195227994156, 1.024619394198223, 1.024619394198223, 1.009765813781591, 0.9744871529235692, 0.9744871529235692, 0.9731820754667071, 0.9911979000151204, 0.9911979000151204, 0.9995532735608961, 0.9995532735608961, 0.9967428578468516)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 3 + $i).Value = $row5[$i]
}

# ---------------------------------------------------------------------
# 4) Rename the "Thomas Hex" method label (row 9's own label - rows keep
#    their own label regardless of the C:T ripple above) to
#    "Matthies Hex".
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 2).Value = "Matthies Hex"
